$d = $word.ActiveDocument

# --- p12 ---
$old = "- Desenvolver um projeto interdisciplinar, de média complexidade, sobre tema relacionado à Engenharia de Produção, similar a situações que os estudantes irão encontrar na vida real, no efetivo exercício de sua profissão; - Aplicar e integrar conhecimentos adquiridos nas demais disciplinas do curso; - Desenvolver competências técnicas (relacionadas ao projeto em si), e competências transversais (aprendizagem ativa, pensamento sistêmico, capacidade de resolução de problemas, trabalho em equipe, liderança, relacionamento interpessoal, gestão de conflitos, capacidade de comunicação, capacidade de planejamento, criatividade e iniciativa) num ambiente de aprendizagem baseado em PBL (Project-Based Learning e Problem-Baed Learning)."
$new = "- Desenvolver um projeto interdisciplinar, de média complexidade, sobre tema relacionado à Engenharia de Produção, similar a situações que os estudantes irão encontrar na vida real, no efetivo exercício de sua profissão; ^l- Aplicar e integrar conhecimentos adquiridos nas demais disciplinas do curso; ^l- Desenvolver competências técnicas (relacionadas ao projeto em si), e competências transversais (aprendizagem ativa, pensamento sistêmico, capacidade de resolução de problemas, trabalho em equipe, liderança, relacionamento interpessoal, gestão de conflitos, capacidade de comunicação, capacidade de planejamento, criatividade e iniciativa) num ambiente de aprendizagem baseado em PBL (Project-Based Learning e Problem-Baed Learning)."
$find = $d.Content.Find
$ok = $find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $ok) { Write-Host "WARNING: p12 replace failed" } else { Write-Host "p12 replaced OK" }

# --- p15 ---
$old = "- Noções de Gestão de Projetos;- Noções de Aprendizagem Baseada em Projetos/Projetos;- Organização do tempo: dimensão pessoal;- Técnicas para a realização de apresentações;- Trabalho em Equipe; - Postura e Ética Profissional;- Técnicas para redação de relatório técnico;- Tutoria de projetos;- Assuntos Técnicos específicos relacionados com o tema do projeto;- Aplicar conhecimentos de Introdução à Engenharia de Produção e Administração e Organização I, Administração e Organização II, Sistemas Produtivos, Estatística, Estatística Multivariada, Economia Geral, Gestão Projetos, Engenharia da Qualidade e Lógica Computacional, integrando-os às demais disciplinas do curso;- Visita (viagem didática complementar) à empresa em que o projeto estiver sendo realizado, para melhor compreender a situação-problema e desenvolver o projeto."
$new = "- Noções de Gestão de Projetos;^l- Noções de Aprendizagem Baseada em Projetos/Projetos;^l- Organização do tempo: dimensão pessoal;^l- Técnicas para a realização de apresentações;^l- Trabalho em Equipe; ^l- Postura e Ética Profissional;^l- Técnicas para redação de relatório técnico;^l- Tutoria de projetos;^l- Assuntos Técnicos específicos relacionados com o tema do projeto;^l- Aplicar conhecimentos de Introdução à Engenharia de Produção e Administração e Organização I, Administração e Organização II, Sistemas Produtivos, Estatística, Estatística Multivariada, Economia Geral, Gestão Projetos, Engenharia da Qualidade e Lógica Computacional, integrando-os às demais disciplinas do curso;^l- Visita (viagem didática complementar) à empresa em que o projeto estiver sendo realizado, para melhor compreender a situação-problema e desenvolver o projeto."
$find = $d.Content.Find
$ok = $find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $ok) { Write-Host "WARNING: p15 replace failed" } else { Write-Host "p15 replaced OK" }

# --- p16 ---
$old = "- O método utilizado tem por fundamento o PBL, que visa desenvolver as competências técnicas relativas ao tema do projeto, bem como competências transversais.- Os estudantes serão divididos em grupos que desenvolverão um projeto durante o semestre relacionado a um tema de Engenharia de Produção, similar ao que eles irão encontrar na vida real no efetivo exercício de sua profissão.- É priorizado o desenvolvimento de projetos em empresas reais.- Cada grupo deverá buscar o conhecimento teórico e prático necessário para ser aplicado no desenvolvimento do projeto.- As aulas ocorrerão: 1) através de uma reunião da equipe de trabalho para tratar do projeto, e 2) palestras e dinâmicas relativas ao tema do projeto, que serão conduzidas por professores ou profissionais de empresas."
$new = "- O método utilizado tem por fundamento o PBL, que visa desenvolver as competências técnicas relativas ao tema do projeto, bem como competências transversais.^l- Os estudantes serão divididos em grupos que desenvolverão um projeto durante o semestre relacionado a um tema de Engenharia de Produção, similar ao que eles irão encontrar na vida real no efetivo exercício de sua profissão.^l- É priorizado o desenvolvimento de projetos em empresas reais.^l- Cada grupo deverá buscar o conhecimento teórico e prático necessário para ser aplicado no desenvolvimento do projeto.^l- As aulas ocorrerão: 1) através de uma reunião da equipe de trabalho para tratar do projeto, e 2) palestras e dinâmicas relativas ao tema do projeto, que serão conduzidas por professores ou profissionais de empresas."
$find = $d.Content.Find
$ok = $find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $ok) { Write-Host "WARNING: p16 replace failed" } else { Write-Host "p16 replaced OK" }

# --- p23 ---
$old = "Notions of Project Management; Notions of Project/Project-Based Learning; Time organization: personal dimension; Techniques for making presentations; Group work, Teams and teams; Professional Posture and Ethics; Techniques for writing a technical report; Project mentoring;Specific technical issues related to the project theme; Apply knowledge of Introduction to Production Engineering and Administration and Organization I, Administration and Organization II, Production Systems, Statistics, Multivariate Statistics, General Economics, Project Management, Quality Engineering and Computational Logic, integrating them with the other subjects of the program; Visit (complementary educational trip) to the company where the project is being carried out, to better understand the problem situation and develop the project."
$new = "Notions of Project Management; Notions of Project/Project-Based Learning; Time organization: personal dimension; Techniques for making presentations; Group work, Teams and teams; Professional Posture and Ethics; Techniques for writing a technical report; Project mentoring;^lSpecific technical issues related to the project theme; Apply knowledge of Introduction to Production Engineering and Administration and Organization I, Administration and Organization II, Production Systems, Statistics, Multivariate Statistics, General Economics, Project Management, Quality Engineering and Computational Logic, integrating them with the other subjects of the program; Visit (complementary educational trip) to the company where the project is being carried out, to better understand the problem situation and develop the project."
$find = $d.Content.Find
$ok = $find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $ok) { Write-Host "WARNING: p23 replace failed" } else { Write-Host "p23 replaced OK" }

# --- p28 ---
$old = "A atividade é denominada Projeto Integrado em Engenharia de Produção II (PIEP) e tem os seguintes objetivos: - Desenvolver projetos ou propostas de solução de problemas reais ou potenciais, de natureza interdisciplinar, de média complexidade, em processos produtivos de empresas, de maneira a contribuir para a capacitação organizacional e o desenvolvimento econômico da região;- Discutir as principais alternativas de projetos ou propostas de solução de problemas com funcionários das empresas participantes da atividade, contribuindo com a capacitação de mão-de-obra das empresas.- Desenvolver competências técnicas e transdisciplinares nos estudantes, entre elas: capacidade de planejamento, desenvolvimento e controle de um projeto; pensamento sistêmico, trabalho em equipe; liderança, capacidade de resolução de problemas, relacionamento interpessoal, gestão de conflitos; argumentação; capacidade de comunicação escrita e falada; ponderação; avaliação crítica e capacidade de tomada de decisão; criatividade e iniciativa; conscientização sobre questões relacionadas com sustentabilidade; construção de protótipos de produtos/sistemas de produção; e, avaliação de proposta de projeto considerando critérios estabelecidos, incluindo critérios econômico-financeiros.- Aplicar os conteúdos das unidades curriculares do curso de Engenharia de Produção no contexto do projeto."
$new = "A atividade é denominada Projeto Integrado em Engenharia de Produção II (PIEP) e tem os seguintes objetivos: ^l- Desenvolver projetos ou propostas de solução de problemas reais ou potenciais, de natureza interdisciplinar, de média complexidade, em processos produtivos de empresas, de maneira a contribuir para a capacitação organizacional e o desenvolvimento econômico da região;^l- Discutir as principais alternativas de projetos ou propostas de solução de problemas com funcionários das empresas participantes da atividade, contribuindo com a capacitação de mão-de-obra das empresas.^l- Desenvolver competências técnicas e transdisciplinares nos estudantes, entre elas: capacidade de planejamento, desenvolvimento e controle de um projeto; pensamento sistêmico, trabalho em equipe; liderança, capacidade de resolução de problemas, relacionamento interpessoal, gestão de conflitos; argumentação; capacidade de comunicação escrita e falada; ponderação; avaliação crítica e capacidade de tomada de decisão; criatividade e iniciativa; conscientização sobre questões relacionadas com sustentabilidade; construção de protótipos de produtos/sistemas de produção; e, avaliação de proposta de projeto considerando critérios estabelecidos, incluindo critérios econômico-financeiros.^l- Aplicar os conteúdos das unidades curriculares do curso de Engenharia de Produção no contexto do projeto."
$find = $d.Content.Find
$ok = $find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $ok) { Write-Host "WARNING: p28 replace failed" } else { Write-Host "p28 replaced OK" }

# --- p30 ---
$old = "A atividade consiste na identificação de uma situação-problema em uma empresa e na tratativa desta situação-problema, de maneira a propor uma solução para a empresa, de maneira dialogada com os profissionais da empresa, contribuindo, também, com a capacitação da mão-de-obra da empresa.Principais etapas da atividade:1.Visita à empresa para identificação da situação-Problema proposta.2.Identificação da teoria relacionada à situação-problema.3.Reuniões com os funcionários da empresa para a identificação das causas reais ou potenciais da situação-problema.4.Identificação das possíveis alternativas para solução da situação-problema.5.Reuniões com os funcionários da empresa para discutir e avaliar as possíveis alternativas para solução da situação-problema.6.Desenvolvimento da Proposta de solução da situação-problema.7.Apresentação e discussão da proposta para os funcionários da empresa8.Realização da avaliação do projeto pela empresa, autoavaliação pelos estudantes e lições aprendidas."
$new = "A atividade consiste na identificação de uma situação-problema em uma empresa e na tratativa desta situação-problema, de maneira a propor uma solução para a empresa, de maneira dialogada com os profissionais da empresa, contribuindo, também, com a capacitação da mão-de-obra da empresa.^lPrincipais etapas da atividade:^l1.Visita à empresa para identificação da situação-Problema proposta.^l2.Identificação da teoria relacionada à situação-problema.^l3.Reuniões com os funcionários da empresa para a identificação das causas reais ou potenciais da situação-problema.^l4.Identificação das possíveis alternativas para solução da situação-problema.^l5.Reuniões com os funcionários da empresa para discutir e avaliar as possíveis alternativas para solução da situação-problema.^l6.Desenvolvimento da Proposta de solução da situação-problema.^l7.Apresentação e discussão da proposta para os funcionários da empresa^l8.Realização da avaliação do projeto pela empresa, autoavaliação pelos estudantes e lições aprendidas."
$find = $d.Content.Find
$ok = $find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $ok) { Write-Host "WARNING: p30 replace failed" } else { Write-Host "p30 replaced OK" }
